$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency list: updated prices and 1h volume
# percentages, plus two rank swaps (rows 19/20 and rows 50/51)
# where Uniswap/TRON and THORChain/FirstDigitalUSD traded places.
#
# Price strings that look numeric (e.g. "0.999", "481.10",
# "0.0000279") get a leading apostrophe so Excel keeps them as
# literal text instead of parsing/rounding them as numbers --
# this matches how the source sheet stores every Price/Volume
# cell as a plain string.

$ws.Range("D2").Value = '67.187.20'
$ws.Range("E2").Value = '  +3.40%  '
$ws.Range("D3").Value = '3.463.37'
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''582.64'
$ws.Range("E5").Value = '  +4.89%  '
$ws.Range("D6").Value = '''190.12'
$ws.Range("E6").Value = '  +8.90%  '
$ws.Range("D7").Value = '''0.633'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '3.453.65'
$ws.Range("E8").Value = '  +2.76%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").Value = '''0.649'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("D12").Value = '''57.36'
$ws.Range("E12").Value = '  +6.64%  '
$ws.Range("D13").Value = '''0.0000279'
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = '''9.51'
$ws.Range("E14").Value = '  +3.70%  '
$ws.Range("D15").Value = '3.998.70'
$ws.Range("E15").Value = '  +2.44%  '
$ws.Range("D16").Value = '''18.92'
$ws.Range("E16").Value = '  +3.50%  '
$ws.Range("D17").Value = '3.449.03'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").Value = '67.106.70'
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''12.18'
$ws.Range("E19").Value = '  +2.28%  '
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").Value = '''0.119'
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("D22").Value = '''481.10'
$ws.Range("E22").Value = '  +5.92%  '
$ws.Range("D23").Value = '''17.11'
$ws.Range("E23").Value = '  +23.05%  '
$ws.Range("D24").Value = '''5.35'
$ws.Range("E24").Value = '  +8.78%  '
$ws.Range("D25").Value = '''4.38'
$ws.Range("E25").Value = '  +7.88%  '
$ws.Range("D26").Value = '''89.77'
$ws.Range("E26").Value = '  +2.59%  '
$ws.Range("E27").Value = '  +4.04%  '
$ws.Range("D28").Value = '''11.02'
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("D29").Value = '''9.12'
$ws.Range("E29").Value = '  +4.93%  '
$ws.Range("D30").Value = '''31.41'
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").Value = '''7.47'
$ws.Range("E31").Value = '  +14.00%  '
$ws.Range("D32").Value = '''11.86'
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").Value = '''599.13'
$ws.Range("E33").Value = '  +3.92%  '
$ws.Range("D34").Value = '''64.28'
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("E35").Value = '  +4.60%  '
$ws.Range("D36").Value = '''0.150'
$ws.Range("E36").Value = '  +6.60%  '
$ws.Range("D38").Value = '''37.54'
$ws.Range("E38").Value = '  +5.13%  '
$ws.Range("D39").Value = '''0.390'
$ws.Range("E39").Value = '  +4.68%  '
$ws.Range("E40").Value = '  -4.44%  '
$ws.Range("D41").Value = '0.0₃0759'
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("D42").Value = '3.240.76'
$ws.Range("E42").Value = '  +5.29%  '
$ws.Range("D43").Value = '''2.96'
$ws.Range("E43").Value = '  +7.20%  '
$ws.Range("D44").Value = '''0.0434'
$ws.Range("E44").Value = '  +4.49%  '
$ws.Range("D45").Value = '''2.88'
$ws.Range("E45").Value = '  +27.65%  '
$ws.Range("D46").Value = '''2.56'
$ws.Range("E46").Value = '  +4.16%  '
$ws.Range("D47").Value = '''3.23'
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("D49").Value = '''3.27'
$ws.Range("E49").Value = '  +12.28%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '''8.73'
$ws.Range("E50").Value = '  +5.55%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = '''0.999'
$ws.Range("E51").Value = '  +0.04%  '
